$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2, 6).Value = 252
$ws1.Cells.Item(3, 6).Value = 446
$ws1.Cells.Item(6, 6).Value = 549
$ws1.Cells.Item(9, 6).Value = 272
$ws1.Cells.Item(11, 6).Value = 352
$ws1.Cells.Item(12, 6).Value = 660
$ws1.Cells.Item(13, 6).Value = 753
$ws1.Cells.Item(14, 6).Value = 1510
$ws1.Cells.Item(15, 6).Value = 1510
$ws1.Cells.Item(16, 6).Value = 884
$ws1.Cells.Item(19, 6).Value = 160
$ws1.Cells.Item(20, 6).Value = 311
$ws1.Cells.Item(23, 6).Value = 101
$ws1.Cells.Item(24, 6).Value = 6581
$ws1.Cells.Item(25, 6).Value = 4933
$ws1.Cells.Item(28, 6).Value = 206
$ws1.Cells.Item(29, 6).Value = 170
$ws1.Cells.Item(32, 6).Value = 1279
$ws1.Cells.Item(34, 6).Value = 248
$ws1.Cells.Item(35, 6).Value = 610
$ws1.Cells.Item(37, 6).Value = 1336
$ws1.Cells.Item(38, 6).Value = 245
$ws1.Cells.Item(40, 6).Value = 145

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(9, 6).Value = 1
$ws2.Cells.Item(18, 6).Value = 238

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(3, 6).Value = 2452
$ws3.Cells.Item(4, 6).Value = 193
$ws3.Cells.Item(5, 6).Value = 53

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(2, 6).Value = 252
$ws4.Cells.Item(4, 6).Value = 446
$ws4.Cells.Item(6, 6).Value = 193
$ws4.Cells.Item(7, 6).Value = 53
$ws4.Cells.Item(9, 6).Value = 549
$ws4.Cells.Item(12, 6).Value = 272
$ws4.Cells.Item(15, 6).Value = 352
$ws4.Cells.Item(16, 6).Value = 660
$ws4.Cells.Item(17, 6).Value = 753
$ws4.Cells.Item(18, 6).Value = 1510
$ws4.Cells.Item(19, 6).Value = 1510
$ws4.Cells.Item(20, 6).Value = 884
$ws4.Cells.Item(23, 6).Value = 160
$ws4.Cells.Item(24, 6).Value = 312
$ws4.Cells.Item(26, 6).Value = 101
$ws4.Cells.Item(29, 6).Value = 6581
$ws4.Cells.Item(30, 6).Value = 4933
$ws4.Cells.Item(33, 6).Value = 1279
$ws4.Cells.Item(35, 6).Value = 248
$ws4.Cells.Item(38, 6).Value = 610
$ws4.Cells.Item(42, 6).Value = 1336
$ws4.Cells.Item(43, 6).Value = 245
$ws4.Cells.Item(44, 6).Value = 145
$ws4.Cells.Item(49, 6).Value = 238
